$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert 3 new rows before row 14, shifting existing rows 14:25 down to 17:28
$ws.Rows("14:16").Insert()

# New block inserted at rows 12 (above the shift point, which is now still at original location)
$ws.Range("A12").Value = "yyy"
$ws.Range("B12").Value = "yyy"
$ws.Range("C12").Value = "noch nicht geliefert"
$ws.Range("D12").Value = "noch nicht relevant"

# New block appended after the shifted table at row 29
$ws.Range("A29").Value = "yyy"
$ws.Range("B29").Value = "yyy"
$ws.Range("C29").Value = "noch nicht geliefert"
$ws.Range("D29").Value = "noch nicht relevant"

$ws.Range("G34").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
